# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.67 = 18462.62 pesos`n✅ 18462.62 pesos = 4.66 = 946.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the transfi rate table ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 214
$ws2.Range("O10").Value = 3951
$ws2.Range("N12").Value = 3957.74
$ws2.Range("O12").Value = 203
